$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("이적 현황")
$ws.Select()
Write-Output $ws.Name
